$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (values that are not ambiguous with numbers)
$ws.Range("D2").Value = "69.125.86"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "2.475.21"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +2.22%  "
$ws.Range("D9").Value = "2.472.85"
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("E10").Value = "  +6.41%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").Value = "69.061.00"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("E20").Value = "  +1.98%  "
$ws.Range("E21").Value = "  +3.24%  "
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("E25").Value = "  +2.88%  "
$ws.Range("D26").Value = "0.0₃0831"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("E27").Value = "  +2.07%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -0.92%  "
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("E45").Value = "  +2.10%  "
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("E50").Value = "  -4.02%  "
$ws.Range("E51").Value = "  -1.57%  "

# Updates whose new text looks like a plain number: force the cell to stay
# text (matching the original inlineStr string cells) by briefly applying a
# text number format, then clearing formatting again so the cell keeps the
# workbook default (unstyled) appearance, exactly like the source cells.
$numericTextCells = @{
    "D5" = "562.48"
    "D6" = "164.45"
    "D8" = "0.512"
    "D10" = "0.158"
    "D12" = "0.334"
    "D13" = "4.86"
    "D15" = "0.0000172"
    "D16" = "23.73"
    "D17" = "10.67"
    "D18" = "339.87"
    "D19" = "6.96"
    "D20" = "3.82"
    "D23" = "66.95"
    "D24" = "3.69"
    "D25" = "8.27"
    "D29" = "432.17"
    "D31" = "1.64"
    "D32" = "158.41"
    "D33" = "19.04"
    "D35" = "0.108"
    "D36" = "17.94"
    "D39" = "1.49"
    "D43" = "131.50"
    "D47" = "0.0920"
    "D51" = "16.95"
}
foreach ($addr in $numericTextCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextCells[$addr]
    $cell.ClearFormats()
}
